$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 191, shifting existing rows 191-271 down to 192-272
$ws.Rows.Item(191).Insert()

# Populate the newly inserted row 191 with the new data
$ws.Cells.Item(191, 1).Value = 11
$ws.Cells.Item(191, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(191, 3).Value = "Bíobío"
$ws.Cells.Item(191, 4).Value = 45027
$ws.Cells.Item(191, 5).Value = 8
$ws.Cells.Item(191, 6).Value = 100112003
$ws.Cells.Item(191, 7).Value = "Ajo"
$ws.Cells.Item(191, 8).Value = "Chino"
$ws.Cells.Item(191, 9).Value = "Primera"
$ws.Cells.Item(191, 10).Value = 250
$ws.Cells.Item(191, 11).Value = 14000
$ws.Cells.Item(191, 12).Value = 15000
$ws.Cells.Item(191, 13).Value = 14400
$ws.Cells.Item(191, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(191, 15).Value = "China"
$ws.Cells.Item(191, 16).Value = 1440
$ws.Cells.Item(191, 17).Value = 10
$ws.Cells.Item(191, 18).Value = "Hortaliza"
